$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add ci.lower / ci.upper columns
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"

# Data rows: ci.lower / ci.upper values
$ws.Range("G2").Value = -0.625348612974379
$ws.Range("H2").Value = -0.106218352586424
$ws.Range("G3").Value = -0.100892665746416
$ws.Range("H3").Value = -0.0230250536387044
$ws.Range("G4").Value = -0.0611683170286652
$ws.Range("H4").Value = 0.0240756422657306
$ws.Range("G5").Value = -0.658551965727739
$ws.Range("H5").Value = -0.150290352847517
$ws.Range("G6").Value = -0.41553321393007
$ws.Range("H6").Value = 0.163552464643116
$ws.Range("G7").Value = -0.67170380173025
$ws.Range("H7").Value = -0.153291777452219
$ws.Range("G8").Value = -0.407234666560044
$ws.Range("H8").Value = 0.160286184524401
$ws.Range("G9").Value = -0.0780617923044647
$ws.Range("H9").Value = -0.00244340476956271
$ws.Range("G10").Value = -0.0585771753447614
$ws.Range("H10").Value = 0.00311086746078643
$ws.Range("G11").Value = -0.0742186618340542
$ws.Range("H11").Value = -0.0126063827881316
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Style = "Normal"
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Style = "Normal"
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Style = "Normal"
$ws.Range("G15").Value = -0.0485806962783365
$ws.Range("H15").Value = 0.0549005438264948
$ws.Range("G16").Value = -0.32343122405823
$ws.Range("H16").Value = 0.365506290596002
$ws.Range("G17").Value = -0.340941539693437
$ws.Range("H17").Value = 0.367391575050634
$ws.Range("G18").Value = -0.12233663405733
$ws.Range("H18").Value = 0.0481512845314613
$ws.Range("G19").Value = -0.814469333120087
$ws.Range("H19").Value = 0.320572369048803
$ws.Range("G20").Value = -0.83106642786014
$ws.Range("H20").Value = 0.327104929286232
$ws.Range("G21").Value = 0.0126063827881316
$ws.Range("H21").Value = 0.0742186618340542
$ws.Range("G22").Value = 0.00488680953912542
$ws.Range("H22").Value = 0.156123584608929
$ws.Range("G23").Value = 0.0839283527916455
$ws.Range("H23").Value = 0.494118744355181
$ws.Range("G24").Value = 0.0325344614643488
$ws.Range("H24").Value = 1.03940959975376
$ws.Range("G25").Value = 0.0717280003454956
$ws.Range("H25").Value = 0.485133568942806
$ws.Range("G26").Value = 0.0262269938164166
$ws.Range("H26").Value = 1.03459607404579
